$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.087.71'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.25%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.939.33'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.10%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '376.15'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.34%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.42'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.00%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.06%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.583'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.88%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.66'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.87%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.04%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0836'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.28%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.402.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.36%  '

# Row 14
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.20%  '

# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.66%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.936.75'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.29%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.20%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.032.36'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.51%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.92%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.14'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.22%  '

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.80%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0950'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.10%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '262.86'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.64%  '

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.23%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.81%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.16'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +9.11%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.74'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +7.03%  '

# Row 28
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.03%  '

# Row 29
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.113'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.71%  '

# Row 30
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.166'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.22%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.65'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.13%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.21%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.08'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.57%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.68'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.61%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0454'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.23%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.16%  '

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.09%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.45%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.57'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.85%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.37%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.114'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.90%  '

# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.65%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '121.60'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.19%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.75%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.80%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.271'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.36%  '

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.28%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.004.52'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.82%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.22'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.08%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0346'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.68%  '

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.94%  '
